$wb = $excel.ActiveWorkbook

# --- Sheet "Sam" (sheet1): delete row 2 (year 2023), shifting all following rows up ---
$wsSam = $wb.Worksheets.Item("Sam")
$wsSam.Activate()
$wsSam.Rows.Item(2).Delete()
$wsSam.Range("A2:XFD2").Select()

# --- Sheet "Casey" (sheet2): delete row 2 (year 2023), shifting all following rows up ---
$wsCasey = $wb.Worksheets.Item("Casey")
$wsCasey.Activate()
$wsCasey.Rows.Item(2).Delete()
$wsCasey.Range("A2:XFD2").Select()

# Leave "Sam" as the active sheet/tab, matching the final saved state
$wsSam.Activate()
